$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D9: "固定基地台" -> "固定基地台，可支援sim卡" (matches C9's content)
$ws.Range("D9").Value = "固定基地台，可支援sim卡"

# Update the active selection to D9 (was E18)
$ws.Range("D9").Select()
